# Update cryptos list rows 2-51 (Coin, Link, Price, Volume(1h))
# Reflects latest scrape: price/volume refresh, a few coins re-ranked,
# and BabyDogeCoin newly entering the top list (TheSandbox drops off).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,('Bitcoin','https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc','29.199.71','  -0.07%  ')
    ,('Ethereum','https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth','1.848.10','  -0.41%  ')
    ,('TetherUSD','https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt','''0.9995','  -0.03%  ')
    ,('BNB','https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb','''245.74','  +1.75%  ')
    ,('XRP','https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp','''0.7028','  +0.88%  ')
    ,('USDC','https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc','''0.9999','  +0.01%  ')
    ,('Dogecoin','https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge','''0.07741','  -0.32%  ')
    ,('Cardano','https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada','''0.3068','  -0.12%  ')
    ,('Solana','https://coinranking.com/coin/zNZHO_Sjf+solana-sol','''23.60','  -0.84%  ')
    ,('TRON','https://coinranking.com/coin/qUhEFk1I61atv+tron-trx','''0.07811','  -0.03%  ')
    ,('Litecoin','https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc','''93.19','  +1.18%  ')
    ,('WrappedEther','https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth','1.853.79','  -0.10%  ')
    ,('Polkadot','https://coinranking.com/coin/25W7FG7om+polkadot-dot','''5.142','  +0.70%  ')
    ,('Polygon','https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic','''0.6867','  -0.15%  ')
    ,('Uniswap','https://coinranking.com/coin/_H5FVG9iW+uniswap-uni','''6.594','  +1.15%  ')
    ,('ShibaInu','https://coinranking.com/coin/xz24e0BjL+shibainu-shib','''0.000008323','  -1.33%  ')
    ,('WrappedBTC','https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc','29.190.78','  -0.10%  ')
    ,('BitcoinCash','https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch','''242.19','  -2.69%  ')
    ,('WrappedliquidstakedEther2.0','https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth','2.091.89','  -0.75%  ')
    ,('Avalanche','https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax','''12.74','  -0.63%  ')
    ,('Dai','https://coinranking.com/coin/MoTuySvg7+dai-dai','''0.9998','  -0.02%  ')
    ,('Chainlink','https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link','''7.508','  -0.18%  ')
    ,('BinanceUSD','https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd','''1.0000','  +0.00%  ')
    ,('Stellar','https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm','''0.1511','  +1.25%  ')
    ,('Monero','https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr','''159.32','  -1.20%  ')
    ,('Cosmos','https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom','''8.846','  -0.28%  ')
    ,('EthereumClassic','https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc','''18.31','  -1.04%  ')
    ,('PancakeSwap','https://coinranking.com/coin/ncYFcP709+pancakeswap-cake','''1.536','  -1.56%  ')
    ,('Filecoin','https://coinranking.com/coin/ymQub4fuB+filecoin-fil','''4.226','  -0.41%  ')
    ,('InternetComputer(DFINITY)','https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp','''4.178','  -0.68%  ')
    ,('Toncoin','https://coinranking.com/coin/67YlI0K1b+toncoin-ton','''1.198','  +0.21%  ')
    ,('Hedera','https://coinranking.com/coin/jad286TjB+hedera-hbar','''0.05120','  -1.60%  ')
    ,('ImmutableX','https://coinranking.com/coin/Z96jIvLU7+immutablex-imx','''0.7899','  +4.04%  ')
    ,('LidoDAOToken','https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo','''1.895','  +2.89%  ')
    ,('ARBITRUM','https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb','''1.148','  -1.80%  ')
    ,('HuobiToken','https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht','''2.696','  -0.39%  ')
    ,('Maker','https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr','1.316.20','  +7.74%  ')
    ,('VeChain','https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet','''0.01867','  +0.23%  ')
    ,('MXToken','https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx','''2.710','  -0.52%  ')
    ,('TrustWalletToken','https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt','''0.9587','  +6.46%  ')
    ,('FraxShare','https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs','''6.078','  +10.28%  ')
    ,('Quant','https://coinranking.com/coin/bauj_21eYVwso+quant-qnt','''106.97','  -2.72%  ')
    ,('PaxDollar','https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp','''0.9999','  +0.09%  ')
    ,('EnergySwap','https://coinranking.com/coin/SbWqqTui-+energyswap-ens','''9.701','  +1.31%  ')
    ,('BabyDogeCoin','https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge','''0.00000000123','  -0.86%  ')
    ,('RocketPoolETH','https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth','1.991.53','  -0.75%  ')
    ,('Mantle','https://coinranking.com/coin/BoI4ux0nd+mantle-mnt','''0.5182','  +0.05%  ')
    ,('Aave','https://coinranking.com/coin/ixgUfzmLR+aave-aave','''64.48','  -1.48%  ')
    ,('RenderToken','https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr','''1.763','  +0.57%  ')
    ,('Aptos','https://coinranking.com/coin/HGYj5JCv5+aptos-apt','''6.982','  -0.87%  ')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
}

Write-Host "Done updating rows 2-51"
